# Registree stats backup on Wed 14 Apr 2021 14:25:22 SAST
#
# 1) Refresh the "as of <date>" timestamp in the title cell (A1) of every
#    worksheet, from 13/04/2021 22:42 -> 14/04/2021 14:25.
# 2) On "410E Attendance", insert a new registree row (Toye, Omolayo, The
#    Wilds, Yes, No) above the old row 93 (Tuckett, Alistair, ...), pushing
#    every following record down by one row, and bump the trailing
#    "Number of attendees" count from 117 to 118.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "as of" timestamp on every sheet's title row ---------
foreach ($s in $wb.Worksheets) {
    $title = $s.Range("A1").Value()
    if ($title) {
        $newTitle = $title.Replace("13/04/2021 22:42", "14/04/2021 14:25")
        if ($newTitle -ne $title) {
            $s.Range("A1").Value = $newTitle
        }
    }
}

# --- 2. Insert the new registree on "410E Attendance" -------------------
$ws = $wb.Worksheets.Item("410E Attendance")

# Push rows 93.. down by one, carrying formatting from the old row 93.
$ws.Rows.Item(93).Insert()
$ws.Range("A94:E94").Copy()
$ws.Range("A93:E93").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(93).RowHeight = 25

# Fill in the new registree's details.
$ws.Range("A93").Value = "Toye"
$ws.Range("B93").Value = "Omolayo"
$ws.Range("C93").Value = "The Wilds"
$ws.Range("D93").Value = "Yes"
$ws.Range("E93").Value = "No"

# The "Number of attendees" summary row shifted from 120 to 121; bump the count.
$ws.Range("A121").Value = "Number of attendees: 118"
